# Update Supplemental Table 5 (education demographics) so that the age-group
# rows no longer use age-standardized values - per commit message
# "updating suppl table 5 so no age standards on age groups".
#
# All affected cells are stored as text (e.g. "0.6", "0.31") rather than
# numbers, so values are written back as text (leading apostrophe forces
# Excel to keep them as text instead of auto-converting to numeric values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ("Less than Primary", first block)
$ws.Range("B2").Value = "'0.59"
$ws.Range("C2").Value = "'0.56"

# Row 3 ("Primary", first block)
$ws.Range("C3").Value = "'0.32"

# Row 4 ("Secondary", first block)
$ws.Range("D4").Value = "'0.22"
$ws.Range("E4").Value = "'0.25"

# Row 8 ("Primary", second block)
$ws.Range("B8").Value = "'0.31"

# Row 9 ("Secondary", second block)
$ws.Range("D9").Value = "'0.21"
$ws.Range("E9").Value = "'0.24"
